$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("fasta-method-1")

$ws.Range("J11").Value = 0.0024149417877200002
$ws.Range("K11").Value = 0.0024149417877200002
$ws.Range("L11").Value = 12.9921875
$ws.Range("M11").Value = 142.96484375

$ws.Range("J12").Value = 0.48018693924
$ws.Range("K12").Value = 0.48260188102700002
$ws.Range("L12").Value = 31.12109375
$ws.Range("M12").Value = 161.03125

$ws.Range("J13").Value = 0.48099708557100002
$ws.Range("K13").Value = 0.48341202735900002
$ws.Range("L13").Value = 31.12109375
$ws.Range("M13").Value = 161.03125

$ws.Range("J14").Value = 0.00168704986572
$ws.Range("K14").Value = 3016.3362939399999
$ws.Range("L14").Value = 241.16015625
$ws.Range("M14").Value = 1103.29296875

$ws.Range("J15").Value = 18.4668970108
$ws.Range("K15").Value = 3034.80319095
$ws.Range("L15").Value = 241.17578125
$ws.Range("M15").Value = 1103.29296875

$ws.Range("J16").Value = 237.360594034
$ws.Range("K16").Value = 3272.1637849799999
$ws.Range("L16").Value = 243.4140625
$ws.Range("M16").Value = 1103.29296875

$ws.Range("J17").Value = 0.0028939247131299999
$ws.Range("K17").Value = 3272.16667891
$ws.Range("L17").Value = 243.4140625
$ws.Range("M17").Value = 1103.29296875

$ws.Activate()
$ws.Range("L26").Select() | Out-Null
